$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Update existing data rows (2-5) with new codes
$ws.Range("A2").Value = "N79426"
$ws.Range("B2").Value = "H04145"
$ws.Range("C2").Value = "F35309"

$ws.Range("A3").Value = "N15841"
$ws.Range("B3").Value = "H96890"
$ws.Range("C3").Value = "F53433"

$ws.Range("A4").Value = "N75760"
$ws.Range("B4").Value = "H26944"
$ws.Range("C4").Value = "F09486"

$ws.Range("A5").Value = "N69119"
$ws.Range("B5").Value = "H54520"
$ws.Range("C5").Value = "F75957"

# Add new row 6
$ws.Range("A6").Value = "N08908"
$ws.Range("B6").Value = "H56185"
$ws.Range("C6").Value = "F65173"
